# Cronograma de proyecto.xlsx - mark four "JUNIO" activities as completed.
#
# Cells O22:O25 and P26 used to hold the literal text "x" (a handwritten
# "done" mark). The author replaced those marks with the numeric value 1,
# which is what the sheet's conditional formatting / COUNTIF-based progress
# formulas (R13:R26, row 28 "Actividades completadas", row 30 "% avance")
# actually key off. Switching "x" -> 1 flips those rows from "Por Completar"
# to "Completado" and drives the overall progress indicators to 100%.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cronograma")

$ws.Range("O22").Value = 1
$ws.Range("O23").Value = 1
$ws.Range("O24").Value = 1
$ws.Range("O25").Value = 1
$ws.Range("P26").Value = 1

# Reflect where the author's cursor ended up after making the last edit.
$ws.Range("P26").Select()
